$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Objetivos:" row now shows the professor info instead of the long
# objectives paragraph.
$ws.Range("B10").Value = "5840560 - Marco Antonio Carvalho Pereira"
$ws.Range("C10").Value = "5840560 - Marco Antonio Carvalho Pereira"

# The old row 13 (no label, just the professor info in B/C) is removed;
# everything below shifts up one row, carrying its row heights with it.
$ws.Rows("13:13").Delete()

# "Programa resumido:" (now row 13) gets a new value.
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# "Programa:" (now row 15) gets a new value.
$ws.Range("B15").Value = "01/01/2015"
$ws.Range("C15").Value = "01/01/2015"

# "Método:" (now row 18) gets a new value.
$ws.Range("B18").Value = "5840560 - Marco Antonio Carvalho Pereira"
$ws.Range("C18").Value = "5840560 - Marco Antonio Carvalho Pereira"
